# New discretization. Control System approach and validation.
#
# Sheet "Deviations" (3rd sheet) gains an additional discretization row
# (row 3) and the existing row 2's last two columns (F2 "XTin", G2 "Qin")
# get new deviation values. All data on this sheet is stored as literal
# text (numbers-as-text), matching the rest of the workbook's convention,
# so every value below is entered with a leading apostrophe to force a
# text cell instead of a numeric one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Deviations")

# --- update existing row 2 ---
$ws.Range("F2").Value = "'0.9956"
$ws.Range("G2").Value = "'0.9"

# --- add new row 3 ---
$ws.Range("A3").Value = "'2"
$ws.Range("B3").Value = "'1"
$ws.Range("C3").Value = "'1"
$ws.Range("D3").Value = "'1"
$ws.Range("E3").Value = "'1"
$ws.Range("F3").Value = "'1.2"
$ws.Range("G3").Value = "'2"
